$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestSuite1")

# Add "null" text to J2 and J3 (error handling values for the new test-case columns)
$ws.Range("J2").Value = "null"
$ws.Range("J3").Value = "null"

# Update the active selection on the sheet to K3
$ws.Range("K3").Select()
